# Update cryptocurrency price/volume data (and one Polkadot/WrappedBTC row swap)
# per the "Updated cryptos list" GitHub Actions commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is a cell address + the new literal text that must be written.
# The sheet stores every data cell (Coin/Link/Price/Volume) as text, including
# price figures that look numeric (e.g. "7.10", "0.998"), so every write below
# forces a Text number format first and restores the Normal style afterward —
# otherwise Excel's COM layer would silently reinterpret number-like strings
# (losing trailing zeros / significant digits, e.g. "7.10" -> 7.0999999999996).
$updates = @(
    @{ Addr = 'D2'; Value = '63.305.99' },
    @{ Addr = 'E2'; Value = '  +2.91%  ' },
    @{ Addr = 'D3'; Value = '3.046.56' },
    @{ Addr = 'E3'; Value = '  +1.73%  ' },
    @{ Addr = 'D4'; Value = '0.998' },
    @{ Addr = 'E4'; Value = '  -0.16%  ' },
    @{ Addr = 'D5'; Value = '596.03' },
    @{ Addr = 'E5'; Value = '  -0.84%  ' },
    @{ Addr = 'D6'; Value = '155.14' },
    @{ Addr = 'E6'; Value = '  +7.41%  ' },
    @{ Addr = 'E7'; Value = '  -0.10%  ' },
    @{ Addr = 'D8'; Value = '3.047.26' },
    @{ Addr = 'E8'; Value = '  +1.82%  ' },
    @{ Addr = 'D9'; Value = '0.518' },
    @{ Addr = 'E9'; Value = '  -0.44%  ' },
    @{ Addr = 'D10'; Value = '6.83' },
    @{ Addr = 'E10'; Value = '  +13.02%  ' },
    @{ Addr = 'E11'; Value = '  +3.79%  ' },
    @{ Addr = 'D12'; Value = '0.468' },
    @{ Addr = 'E12'; Value = '  +2.44%  ' },
    @{ Addr = 'D13'; Value = '0.0000236' },
    @{ Addr = 'E13'; Value = '  +2.91%  ' },
    @{ Addr = 'D14'; Value = '35.83' },
    @{ Addr = 'E14'; Value = '  +4.15%  ' },
    @{ Addr = 'D16'; Value = '3.545.03' },
    @{ Addr = 'E16'; Value = '  +1.49%  ' },
    @{ Addr = 'B17'; Value = 'Polkadot' },
    @{ Addr = 'C17'; Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot' },
    @{ Addr = 'D17'; Value = '7.10' },
    @{ Addr = 'E17'; Value = '  +1.84%  ' },
    @{ Addr = 'B18'; Value = 'WrappedBTC' },
    @{ Addr = 'C18'; Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc' },
    @{ Addr = 'D18'; Value = '63.130.01' },
    @{ Addr = 'E18'; Value = '  +2.67%  ' },
    @{ Addr = 'D19'; Value = '3.043.13' },
    @{ Addr = 'E19'; Value = '  +1.50%  ' },
    @{ Addr = 'D20'; Value = '456.10' },
    @{ Addr = 'E20'; Value = '  +0.96%  ' },
    @{ Addr = 'D21'; Value = '14.35' },
    @{ Addr = 'E21'; Value = '  +2.37%  ' },
    @{ Addr = 'D22'; Value = '0.701' },
    @{ Addr = 'E22'; Value = '  +2.16%  ' },
    @{ Addr = 'D23'; Value = '7.54' },
    @{ Addr = 'E23'; Value = '  +2.81%  ' },
    @{ Addr = 'D24'; Value = '83.22' },
    @{ Addr = 'E24'; Value = '  +2.01%  ' },
    @{ Addr = 'D25'; Value = '11.34' },
    @{ Addr = 'E25'; Value = '  +4.97%  ' },
    @{ Addr = 'D26'; Value = '2.32' },
    @{ Addr = 'E26'; Value = '  +4.09%  ' },
    @{ Addr = 'D27'; Value = '12.41' },
    @{ Addr = 'E27'; Value = '  +3.74%  ' },
    @{ Addr = 'E28'; Value = '  +0.05%  ' },
    @{ Addr = 'D29'; Value = '7.54' },
    @{ Addr = 'E29'; Value = '  +4.72%  ' },
    @{ Addr = 'E30'; Value = '  +0.38%  ' },
    @{ Addr = 'D31'; Value = '2.25' },
    @{ Addr = 'E31'; Value = '  +8.80%  ' },
    @{ Addr = 'D32'; Value = '0.997' },
    @{ Addr = 'E32'; Value = '  -0.33%  ' },
    @{ Addr = 'D33'; Value = '27.70' },
    @{ Addr = 'E33'; Value = '  +1.27%  ' },
    @{ Addr = 'E34'; Value = '  +0.87%  ' },
    @{ Addr = 'D35'; Value = '0.0₃0867' },
    @{ Addr = 'E35'; Value = '  +4.74%  ' },
    @{ Addr = 'D36'; Value = '1.05' },
    @{ Addr = 'E36'; Value = '  +2.44%  ' },
    @{ Addr = 'D37'; Value = '5.97' },
    @{ Addr = 'E37'; Value = '  +3.02%  ' },
    @{ Addr = 'D38'; Value = '3.24' },
    @{ Addr = 'E38'; Value = '  +12.69%  ' },
    @{ Addr = 'D39'; Value = '2.13' },
    @{ Addr = 'E39'; Value = '  +3.18%  ' },
    @{ Addr = 'D40'; Value = '0.131' },
    @{ Addr = 'E40'; Value = '  +5.57%  ' },
    @{ Addr = 'D41'; Value = '50.52' },
    @{ Addr = 'E41'; Value = '  +0.27%  ' },
    @{ Addr = 'D42'; Value = '9.14' },
    @{ Addr = 'E42'; Value = '  -0.58%  ' },
    @{ Addr = 'D43'; Value = '0.305' },
    @{ Addr = 'E43'; Value = '  +12.59%  ' },
    @{ Addr = 'D44'; Value = '43.62' },
    @{ Addr = 'E44'; Value = '  +9.68%  ' },
    @{ Addr = 'D45'; Value = '397.85' },
    @{ Addr = 'E45'; Value = '  -0.18%  ' },
    @{ Addr = 'D46'; Value = '0.0363' },
    @{ Addr = 'E46'; Value = '  +2.73%  ' },
    @{ Addr = 'D47'; Value = '2.729.13' },
    @{ Addr = 'E47'; Value = '  +1.38%  ' },
    @{ Addr = 'D48'; Value = '132.07' },
    @{ Addr = 'E48'; Value = '  +0.75%  ' },
    @{ Addr = 'D49'; Value = '2.30' },
    @{ Addr = 'E49'; Value = '  +7.09%  ' },
    @{ Addr = 'D51'; Value = '24.58' },
    @{ Addr = 'E51'; Value = '  +4.28%  ' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Addr)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.Style = "Normal"
}
